$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected (password "D382"); unprotect so the cells can be edited.
$ws.Unprotect("D382")

# Update the confidential/disclaimer note: the "as of" date moved from 2021-05-11 to 2021-05-12.
$ws.Range("A12").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-12 for illustrative purposes only and are subject to change."

# Refresh the Weight / Percent Change figures for rows 2-9.
$ws.Range("D2").Value2 = 0.1773447658607328
$ws.Range("E2").Value2 = -0.002705139765554532

$ws.Range("D3").Value2 = 0.1773847443922613
$ws.Range("E3").Value2 = -0.002941176470588225

$ws.Range("D4").Value2 = 0.2255878593195454
$ws.Range("E4").Value2 = -0.005799502899751507

$ws.Range("D5").Value2 = 0.07982313497651761
$ws.Range("E5").Value2 = -0.001014198782961384

$ws.Range("D6").Value2 = 0.07962923909860405

$ws.Range("D7").Value2 = 0.1201774647014553
$ws.Range("E7").Value2 = -0.001962708537782021

$ws.Range("D8").Value2 = 0.1400527916508835
$ws.Range("E8").Value2 = -0.00248344370860909

$ws.Range("E9").Value2 = -0.002974402745725579

# Restore sheet protection with the original password.
$ws.Protect("D382")
